$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (STRASSE, HAUSNR) before the old PLZ column (J),
# shifting PLZ/WOHNORT/VERMOEGEN/EL-BEZUG/SH-BEZUG two columns to the right.
$ws.Range("J1:K1").EntireColumn.Insert() | Out-Null

# Match the new columns' width to the neighbouring AUFENTHALTSBEWILLIGUNG column (I).
$ws.Range("J1:K1").ColumnWidth = $ws.Range("I1").ColumnWidth

# Header row
$ws.Range("J1").Value = "STRASSE"
$ws.Range("K1").Value = "HAUSNR"

# Row 2 (Ackermann / Aaron)
$ws.Range("J2").Value = "Ackerstrasse"
$ws.Range("K2").Value = 11

# Row 4 (Freud / Ferdi) - set before row 3 so shared strings land in the same
# order as the authored workbook (Fichtenhalde before Denzingsteig).
$ws.Range("J4").Value = "Fichtenhalde"
$ws.Range("K4").Value = 300

# Row 3 (Da Silva / Denise)
$ws.Range("J3").Value = "Denzingsteig"
$ws.Range("K3").Value = 22

# Restore the selection recorded in the authored workbook.
$ws.Range("M21").Select() | Out-Null
